# Apply KHL probabilities update (2025-12-02 publish)
$wb = $excel.ActiveWorkbook
$wsSummary = $wb.Worksheets.Item("Summary")
$wsCards = $wb.Worksheets.Item("Cards_telegram")

# --- Summary sheet, row 2 ---
$wsSummary.Range("B2").Value = 45993.51041666666
$wsSummary.Range("C2").Value = "Амур"
$wsSummary.Range("D2").Value = "Динамо Мн"
$wsSummary.Range("E2").Value = "Амур – Динамо Мн"
$wsSummary.Range("F2").Value = 897838
$wsSummary.Range("G2").Value = "https://text.khl.ru/text/897838.html"
$wsSummary.Range("H2").Value = 1.1
$wsSummary.Range("I2").Value = 4.615385
$wsSummary.Range("J2").Value = 5.715385
$wsSummary.Range("K2").Value = 22.192716
$wsSummary.Range("L2").Value = 42.630188
$wsSummary.Range("M2").Value = 64.82290399999999
$wsSummary.Range("N2").Value = 0.110222
$wsSummary.Range("O2").Value = 0.134865
$wsSummary.Range("P2").Value = 0.7545460000000001
$wsSummary.Range("Q2").Value = 9.072598936691405
$wsSummary.Range("R2").Value = 7.414822229637044
$wsSummary.Range("S2").Value = 1.325300246770906
$wsSummary.Range("T2").Value = 11.0222
$wsSummary.Range("U2").Value = 13.4865
$wsSummary.Range("V2").Value = 75.4546
$wsSummary.Range("W2").Value = 0.580246
$wsSummary.Range("X2").Value = 0.419386
$wsSummary.Range("Y2").Value = 2.384438202515105
$wsSummary.Range("Z2").Value = 0.745048
$wsSummary.Range("AA2").Value = 0.254585
$wsSummary.Range("AB2").Value = 3.927961191743425
$wsSummary.Range("AC2").Value = 0.861764
$wsSummary.Range("AD2").Value = 0.137868
$wsSummary.Range("AE2").Value = 7.253314764847536
$wsSummary.Range("AF2").Value = 0.324827
$wsSummary.Range("AG2").Value = 0.675173
$wsSummary.Range("AH2").Value = 3.078561819060608
$wsSummary.Range("AI2").Value = 0.113094
$wsSummary.Range("AJ2").Value = 0.886906
$wsSummary.Range("AK2").Value = 8.842202062001521
$wsSummary.Range("AL2").Value = 0.813045
$wsSummary.Range("AM2").Value = 0.186955
$wsSummary.Range("AN2").Value = 1.229944222029531
$wsSummary.Range("AO2").Value = 0.595354
$wsSummary.Range("AP2").Value = 0.404646
$wsSummary.Range("AQ2").Value = 1.679672934086275
$wsSummary.Range("AR2").Value = 0.432665
$wsSummary.Range("AS2").Value = 2.311256977107
$wsSummary.Range("AT2").Value = 0.960294
$wsSummary.Range("AU2").Value = 1.041347753917029

# --- Summary sheet, row 3 ---
$wsSummary.Range("B3").Value = 45993.52083333334
$wsSummary.Range("C3").Value = "Адмирал"
$wsSummary.Range("D3").Value = "ХК Сочи"
$wsSummary.Range("E3").Value = "Адмирал – ХК Сочи"
$wsSummary.Range("F3").Value = 897839
$wsSummary.Range("G3").Value = "https://text.khl.ru/text/897839.html"
$wsSummary.Range("H3").Value = 2.505511
$wsSummary.Range("I3").Value = 0.928571
$wsSummary.Range("J3").Value = 3.434082
$wsSummary.Range("K3").Value = 32.125309
$wsSummary.Range("L3").Value = 21.640231
$wsSummary.Range("M3").Value = 53.76554
$wsSummary.Range("N3").Value = 0.856225
$wsSummary.Range("O3").Value = 0.07615
$wsSummary.Range("P3").Value = 0.060739
$wsSummary.Range("Q3").Value = 1.167917311454349
$wsSummary.Range("R3").Value = 13.13197636244255
$wsSummary.Range("S3").Value = 16.46388646503894
$wsSummary.Range("T3").Value = 85.6225
$wsSummary.Range("U3").Value = 7.614999999999999
$wsSummary.Range("V3").Value = 6.0739
$wsSummary.Range("W3").Value = 0.292915
$wsSummary.Range("X3").Value = 0.700198
$wsSummary.Range("Y3").Value = 1.428167461203831
$wsSummary.Range("Z3").Value = 0.455064
$wsSummary.Range("AA3").Value = 0.538049
$wsSummary.Range("AB3").Value = 1.858566784809562
$wsSummary.Range("AC3").Value = 0.615642
$wsSummary.Range("AD3").Value = 0.377472
$wsSummary.Range("AE3").Value = 2.649203119701594
$wsSummary.Range("AF3").Value = 0.939923
$wsSummary.Range("AG3").Value = 0.060077
$wsSummary.Range("AH3").Value = 1.063916937876826
$wsSummary.Range("AI3").Value = 0.828727
$wsSummary.Range("AJ3").Value = 0.171273
$wsSummary.Range("AK3").Value = 1.20666998903137
$wsSummary.Range("AL3").Value = 0.415477
$wsSummary.Range("AM3").Value = 0.584523
$wsSummary.Range("AN3").Value = 2.406872101223413
$wsSummary.Range("AO3").Value = 0.171658
$wsSummary.Range("AP3").Value = 0.828342
$wsSummary.Range("AQ3").Value = 5.825536823218259
$wsSummary.Range("AR3").Value = 0.970588
$wsSummary.Range("AS3").Value = 1.030303280073522
$wsSummary.Range("AT3").Value = 0.258435
$wsSummary.Range("AU3").Value = 3.86944492812506

# --- Summary sheet, row 4 ---
$wsSummary.Range("B4").Value = 45993.79166666666
$wsSummary.Range("C4").Value = "Локомотив"
$wsSummary.Range("D4").Value = "СКА"
$wsSummary.Range("E4").Value = "Локомотив – СКА"
$wsSummary.Range("F4").Value = 897840
$wsSummary.Range("G4").Value = "https://text.khl.ru/text/897840.html"
$wsSummary.Range("H4").Value = 2.392003
$wsSummary.Range("I4").Value = 3.857143
$wsSummary.Range("J4").Value = 6.249146
$wsSummary.Range("K4").Value = 28.33909
$wsSummary.Range("L4").Value = 33.415371
$wsSummary.Range("M4").Value = 61.75446
$wsSummary.Range("N4").Value = 0.401901
$wsSummary.Range("O4").Value = 0.162865
$wsSummary.Range("P4").Value = 0.434399
$wsSummary.Range("Q4").Value = 2.488174948556983
$wsSummary.Range("R4").Value = 6.140054646486353
$wsSummary.Range("S4").Value = 2.30203108202367
$wsSummary.Range("T4").Value = 40.1901
$wsSummary.Range("U4").Value = 16.2865
$wsSummary.Range("V4").Value = 43.43989999999999
$wsSummary.Range("W4").Value = 0.251846
$wsSummary.Range("X4").Value = 0.747318
$wsSummary.Range("Y4").Value = 1.338118444892268
$wsSummary.Range("Z4").Value = 0.404979
$wsSummary.Range("AA4").Value = 0.594185
$wsSummary.Range("AB4").Value = 1.682977523835169
$wsSummary.Range("AC4").Value = 0.56473
$wsSummary.Range("AD4").Value = 0.434434
$wsSummary.Range("AE4").Value = 2.30184561981797
$wsSummary.Range("AF4").Value = 0.812122
$wsSummary.Range("AG4").Value = 0.187878
$wsSummary.Range("AH4").Value = 1.231342088996481
$wsSummary.Range("AI4").Value = 0.593931
$wsSummary.Range("AJ4").Value = 0.406069
$wsSummary.Range("AK4").Value = 1.683697264497054
$wsSummary.Range("AL4").Value = 0.826419
$wsSummary.Range("AM4").Value = 0.173581
$wsSummary.Range("AN4").Value = 1.210039943418532
$wsSummary.Range("AO4").Value = 0.616299
$wsSummary.Range("AP4").Value = 0.383701
$wsSummary.Range("AQ4").Value = 1.622589035516851
$wsSummary.Range("AR4").Value = 0.716493
$wsSummary.Range("AS4").Value = 1.395687047884627
$wsSummary.Range("AT4").Value = 0.744011
$wsSummary.Range("AU4").Value = 1.344066149559617

# Remove the now-absent 5th Summary row (match dropped from this slate)
$wsSummary.Rows.Item(5).Delete()

# --- Cards_telegram sheet, row 2 ---
$wsCards.Range("A2").Value = 45993.51041666666
$wsCards.Range("B2").Value = "Амур – Динамо Мн"
$cardText2 = @"
КХЛ • Регулярный чемпионат • 02.12.2025

Амур – Динамо Мн

Ожидания модели (60’):
• Голы: λ_total ≈ 4.25 (1.17 : 3.08)
• Броски: SOG λ ≈ 65 (22 : 43)

Исход (60’), честные кф:
• П1: 11.0%  (Kмод 9.07)
• Х:  13.5%  (Kмод 7.41)
• П2: 75.5%  (Kмод 1.33)

Тоталы голов:
• ТМ 4.5: 58.0%  (Kмод 1.72)
• ТБ 4.5: 41.9%  (Kмод 2.38)

• ТМ 5.5: 74.5%  (Kмод 1.34)
• ТБ 5.5: 25.5%  (Kмод 3.93)

• ТМ 6.5: 86.2%  (Kмод 1.16)
• ТБ 6.5: 13.8%  (Kмод 7.25)

Индивидуальные тоталы:
• Амур ИТБ 1.5: 32.5% (Kмод 3.08)
• Амур ИТБ 2.5: 11.3% (Kмод 8.84)
• Динамо Мн ИТБ 1.5: 81.3% (Kмод 1.23)
• Динамо Мн ИТБ 2.5: 59.5% (Kмод 1.68)

Фора +1.5:
• Амур +1.5: 43.3% (Kмод 2.31)
• Динамо Мн +1.5: 96.0% (Kмод 1.04)
"@
$wsCards.Range("C2").Value = $cardText2

# --- Cards_telegram sheet, row 3 ---
$wsCards.Range("A3").Value = 45993.52083333334
$wsCards.Range("B3").Value = "Адмирал – ХК Сочи"
$cardText3 = @"
КХЛ • Регулярный чемпионат • 02.12.2025

Адмирал – ХК Сочи

Ожидания модели (60’):
• Голы: λ_total ≈ 5.94 (4.52 : 1.42)
• Броски: SOG λ ≈ 54 (32 : 22)

Исход (60’), честные кф:
• П1: 85.6%  (Kмод 1.17)
• Х:  7.6%  (Kмод 13.13)
• П2: 6.1%  (Kмод 16.46)

Тоталы голов:
• ТМ 4.5: 29.3%  (Kмод 3.41)
• ТБ 4.5: 70.0%  (Kмод 1.43)

• ТМ 5.5: 45.5%  (Kмод 2.20)
• ТБ 5.5: 53.8%  (Kмод 1.86)

• ТМ 6.5: 61.6%  (Kмод 1.62)
• ТБ 6.5: 37.7%  (Kмод 2.65)

Индивидуальные тоталы:
• Адмирал ИТБ 1.5: 94.0% (Kмод 1.06)
• Адмирал ИТБ 2.5: 82.9% (Kмод 1.21)
• ХК Сочи ИТБ 1.5: 41.5% (Kмод 2.41)
• ХК Сочи ИТБ 2.5: 17.2% (Kмод 5.83)

Фора +1.5:
• Адмирал +1.5: 97.1% (Kмод 1.03)
• ХК Сочи +1.5: 25.8% (Kмод 3.87)
"@
$wsCards.Range("C3").Value = $cardText3

# --- Cards_telegram sheet, row 4 ---
$wsCards.Range("A4").Value = 45993.79166666666
$wsCards.Range("B4").Value = "Локомотив – СКА"
$cardText4 = @"
КХЛ • Регулярный чемпионат • 02.12.2025

Локомотив – СКА

Ожидания модели (60’):
• Голы: λ_total ≈ 6.26 (3.08 : 3.18)
• Броски: SOG λ ≈ 62 (28 : 33)

Исход (60’), честные кф:
• П1: 40.2%  (Kмод 2.49)
• Х:  16.3%  (Kмод 6.14)
• П2: 43.4%  (Kмод 2.30)

Тоталы голов:
• ТМ 4.5: 25.2%  (Kмод 3.97)
• ТБ 4.5: 74.7%  (Kмод 1.34)

• ТМ 5.5: 40.5%  (Kмод 2.47)
• ТБ 5.5: 59.4%  (Kмод 1.68)

• ТМ 6.5: 56.5%  (Kмод 1.77)
• ТБ 6.5: 43.4%  (Kмод 2.30)

Индивидуальные тоталы:
• Локомотив ИТБ 1.5: 81.2% (Kмод 1.23)
• Локомотив ИТБ 2.5: 59.4% (Kмод 1.68)
• СКА ИТБ 1.5: 82.6% (Kмод 1.21)
• СКА ИТБ 2.5: 61.6% (Kмод 1.62)

Фора +1.5:
• Локомотив +1.5: 71.6% (Kмод 1.40)
• СКА +1.5: 74.4% (Kмод 1.34)
"@
$wsCards.Range("C4").Value = $cardText4

# Remove the now-absent 5th Cards_telegram row (match dropped from this slate)
$wsCards.Rows.Item(5).Delete()